# Week 13 logging update
# - Rushing sheet: new player "M.Boone" inserted as row 5 (between J.Williams and
#   D.Crockett), pushing D.Crockett / J.Jeudy / D.Hamilton down one row each.
#   T.Bridgewater and J.Williams rushing totals also bump up for the week.
# - Receiving sheet: stat totals bump up for several players (no row/order change).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Insert a new row 5 for M.Boone; existing rows 5-7 (D.Crockett, J.Jeudy,
# D.Hamilton) shift down to rows 6-8, formatting and all.
$rushing.Rows(5).Insert()

# The inserted row's first cell doesn't inherit the bordered/bold "index"
# style used by the rest of column A - copy it over from the row above.
$rushing.Range("A4").Copy()
$rushing.Range("A5").PasteSpecial(-4122)

# New M.Boone row.
$rushing.Range("A5").Value = 3
$rushing.Range("B5").Value = "M.Boone"
$rushing.Range("C5").Value = 3
$rushing.Range("D5").Value = 1
$rushing.Range("E5").Value = 0
$rushing.Range("F5").Value = 0

# Renumber the shifted rows' index column.
$rushing.Range("A6").Value = 4
$rushing.Range("A7").Value = 5
$rushing.Range("A8").Value = 6

# Updated weekly totals.
$rushing.Range("D2").Value = 6
$rushing.Range("E2").Value = 11
$rushing.Range("F2").Value = 8

$rushing.Range("C4").Value = 73
$rushing.Range("D4").Value = 51
$rushing.Range("E4").Value = 15
$rushing.Range("F4").Value = 19

# ---------------------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# J.Williams
$receiving.Range("C3").Value = 40
$receiving.Range("D3").Value = 31
$receiving.Range("E3").Value = 2
$receiving.Range("F3").Value = 2
$receiving.Range("G3").Value = 6
$receiving.Range("H3").Value = 4

# M.Boone
$receiving.Range("C4").Value = 2
$receiving.Range("D4").Value = 2

# C.Sutton
$receiving.Range("C5").Value = 51
$receiving.Range("D5").Value = 39
$receiving.Range("E5").Value = 25

# J.Jeudy
$receiving.Range("C6").Value = 65
$receiving.Range("D6").Value = 53
$receiving.Range("E6").Value = 19
$receiving.Range("F6").Value = 15

# T.Patrick
$receiving.Range("C7").Value = 48
$receiving.Range("D7").Value = 33
$receiving.Range("E7").Value = 13

# K.Hinton
$receiving.Range("C9").Value = 4
$receiving.Range("D9").Value = 3

# N.Fant
$receiving.Range("C10").Value = 63
$receiving.Range("D10").Value = 50

# A.Okwuegbunam
$receiving.Range("C11").Value = 25
$receiving.Range("D11").Value = 23
$receiving.Range("G11").Value = 3
